# Update leve-profit calculation cells across sheets (scheduled market-price refresh).
# Values below mirror the authoritative diff; $null clears a cell entirely (removes the <c> element).
$wb = $excel.ActiveWorkbook

# ALC!row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2179618.5
$ws.Range("J17").Value = 2179618.5
$ws.Range("L17").Value = 6538855.5
$ws.Range("N17").Value = -6539191.5

# ALC!row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2902.25
$ws.Range("I40").Value = 2895
$ws.Range("J40").Value = 2909.5
$ws.Range("K40").Value = 2895
$ws.Range("L40").Value = 2909.5
$ws.Range("M40").Value = -2720
$ws.Range("N40").Value = -3259.5

# ALC!row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 16243.813
$ws.Range("J87").Value = 16243.813
$ws.Range("L87").Value = 16243.813
$ws.Range("N87").Value = -18739.813

# ALC!row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 16243.813
$ws.Range("J90").Value = 16243.813
$ws.Range("L90").Value = 48731.439
$ws.Range("N90").Value = -61211.439

# ALC!row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1261.6666
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 1256.875
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 3770.625
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -5986.625

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 30772.883
$ws.Range("I132").Value = 41385.29
$ws.Range("J132").Value = 3357.5
$ws.Range("K132").Value = 124155.87
$ws.Range("L132").Value = 10072.5
$ws.Range("M132").Value = -121625.87
$ws.Range("N132").Value = -15132.5

# ALC!row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 45686
$ws.Range("J134").Value = 45686
$ws.Range("L134").Value = 45686
$ws.Range("N134").Value = -55826

# ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 989.53845
$ws.Range("I135").Value = 1082.1052
$ws.Range("J135").Value = 738.2857
$ws.Range("K135").Value = 9738.9468
$ws.Range("L135").Value = 6644.571300000001
$ws.Range("M135").Value = -7203.9468
$ws.Range("N135").Value = -11714.5713

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 910.875
$ws.Range("I137").Value = 800.6667
$ws.Range("J137").Value = 977
$ws.Range("K137").Value = 2402.0001
$ws.Range("L137").Value = 2931
$ws.Range("M137").Value = 147.9998999999998
$ws.Range("N137").Value = -8031

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1936.35
$ws.Range("I138").Value = 1456.2245
$ws.Range("J138").Value = 2397.647
$ws.Range("K138").Value = 4368.6735
$ws.Range("L138").Value = 7192.941
$ws.Range("M138").Value = 771.3265000000001
$ws.Range("N138").Value = -17472.941

# ARM!row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 9177.5
$ws.Range("J24").Value = 9177.5
$ws.Range("L24").Value = 9177.5
$ws.Range("N24").Value = -9925.5

# ARM!row 28
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2793.8462
$ws.Range("I28").Value = 2793.8462
$ws.Range("K28").Value = 2793.8462
$ws.Range("M28").Value = -2601.8462

# ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1953
$ws.Range("I74").Value = 926
$ws.Range("J74").Value = 3273.4285
$ws.Range("K74").Value = 926
$ws.Range("L74").Value = 3273.4285
$ws.Range("M74").Value = -52
$ws.Range("N74").Value = -5021.4285

# ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1953
$ws.Range("I77").Value = 926
$ws.Range("J77").Value = 3273.4285
$ws.Range("K77").Value = 4630
$ws.Range("L77").Value = 16367.1425
$ws.Range("M77").Value = -262
$ws.Range("N77").Value = -25103.1425

# ARM!row 99
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 2793.8462
$ws.Range("I99").Value = 2793.8462
$ws.Range("K99").Value = 2793.8462
$ws.Range("M99").Value = 201.1538

# ARM!row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 9177.5
$ws.Range("J100").Value = 9177.5
$ws.Range("L100").Value = 9177.5
$ws.Range("N100").Value = -11341.5

# ARM!row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

# ARM!row 107
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 23114
$ws.Range("J107").Value = 23114
$ws.Range("L107").Value = 23114
$ws.Range("N107").Value = -30794

# ARM!row 111
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = $null

# ARM!row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2733.3333
$ws.Range("I122").Value = 1633.3334
$ws.Range("K122").Value = 4900.0002
$ws.Range("M122").Value = -2450.0002

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2066.4524
$ws.Range("I132").Value = 1691.6072
$ws.Range("J132").Value = 2816.1428
$ws.Range("K132").Value = 5074.821599999999
$ws.Range("L132").Value = 8448.428400000001
$ws.Range("M132").Value = -2544.821599999999
$ws.Range("N132").Value = -13508.4284

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1990.7609
$ws.Range("I31").Value = 1273
$ws.Range("K31").Value = 1273
$ws.Range("M31").Value = -978

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1990.7609
$ws.Range("I34").Value = 1273
$ws.Range("K34").Value = 1273
$ws.Range("M34").Value = -1071

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2163.7026
$ws.Range("I58").Value = 1710.4445
$ws.Range("K58").Value = 1710.4445
$ws.Range("M58").Value = -1507.4445

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2163.7026
$ws.Range("I136").Value = 1710.4445
$ws.Range("K136").Value = 5131.333500000001
$ws.Range("M136").Value = -2581.333500000001

# CUL!row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 35.285713
$ws.Range("I2").Value = 81.2
$ws.Range("J2").Value = 9.777778
$ws.Range("K2").Value = 487.2
$ws.Range("L2").Value = 58.666668
$ws.Range("M2").Value = -374.2
$ws.Range("N2").Value = -284.666668

# CUL!row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 125180.125
$ws.Range("I38").Value = 63.333332
$ws.Range("J38").Value = 200250.2
$ws.Range("K38").Value = 189.999996
$ws.Range("L38").Value = 600750.6000000001
$ws.Range("M38").Value = 157.000004
$ws.Range("N38").Value = -601444.6000000001

# CUL!row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 917.44446
$ws.Range("I122").Value = 505.16
$ws.Range("J122").Value = 1854.4546
$ws.Range("K122").Value = 4546.440000000001
$ws.Range("L122").Value = 16690.0914
$ws.Range("M122").Value = -2096.440000000001
$ws.Range("N122").Value = -21590.0914

# CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 932.67
$ws.Range("I131").Value = 543.3333
$ws.Range("J131").Value = 944.71136
$ws.Range("K131").Value = 1629.9999
$ws.Range("L131").Value = 2834.13408
$ws.Range("M131").Value = 3410.0001
$ws.Range("N131").Value = -12914.13408

# CUL!row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6060.1704
$ws.Range("I137").Value = 1894.8572
$ws.Range("J137").Value = 7827.273
$ws.Range("K137").Value = 5684.571599999999
$ws.Range("L137").Value = 23481.819
$ws.Range("M137").Value = -584.5715999999993
$ws.Range("N137").Value = -33681.819

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2347.2258
$ws.Range("I132").Value = 1959.619
$ws.Range("K132").Value = 5878.857
$ws.Range("M132").Value = -3348.857

# LTW!row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2175.15
$ws.Range("I16").Value = 2218.2727
$ws.Range("J16").Value = 2122.4443
$ws.Range("K16").Value = 2218.2727
$ws.Range("L16").Value = 2122.4443
$ws.Range("M16").Value = -2048.2727
$ws.Range("N16").Value = -2462.4443

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5763.8135
$ws.Range("I132").Value = 5817.302
$ws.Range("J132").Value = 5634.9546
$ws.Range("K132").Value = 17451.906
$ws.Range("L132").Value = 16904.8638
$ws.Range("M132").Value = -14921.906
$ws.Range("N132").Value = -21964.8638
